$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$shp = $hdr.Range.InlineShapes.Item(1)
try {
  $xml = $shp.WordOpenXML
  Write-Output ("LEN=" + $xml.Length)
  Write-Output $xml
} catch {
  Write-Output ("err: " + $_.Exception.Message)
}
